$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-10-17 Thursday" "2024-10-18 Friday"

Replace-Text "32×38=" "60×14="
Replace-Text "18×85=" "26×13="
Replace-Text "69×67=" "35×25="
Replace-Text "18×31=" "32×41="
Replace-Text "13×94=" "35×50="
Replace-Text "71×39=" "18×73="
Replace-Text "34×23=" "79×64="
Replace-Text "83×71=" "28×71="
Replace-Text "26×81=" "14×55="
Replace-Text "89×86=" "18×60="
Replace-Text "46×88=" "56×88="
Replace-Text "58×31=" "61×95="
Replace-Text "15×18=" "89×92="
Replace-Text "78×89=" "74×34="
Replace-Text "63×47=" "73×59="
Replace-Text "53×91=" "88×80="
Replace-Text "56×85=" "20×32="
Replace-Text "13×50=" "23×96="
Replace-Text "42×86=" "37×71="
Replace-Text "59×52=" "83×36="
Replace-Text "15×75=" "73×51="
Replace-Text "35×12=" "44×37="
Replace-Text "14×56=" "87×82="
Replace-Text "42×52=" "12×76="
Replace-Text "22×88=" "85×85="
